$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply green fill + keep wrap text for column A (A5:A10) - matches the wrap style used elsewhere
$rangeA = $ws.Range("A5:A10")
$rangeA.Interior.Color = 5296274   # RGB(146,208,80) == hex 92D050
$rangeA.WrapText = $true

# Apply green fill (no wrap) for the rest of the used block, B5:F10 plus G5:G6
$rangeRest = $ws.Range("B5:F10")
$rangeRest.Interior.Color = 5296274

$rangeG = $ws.Range("G5:G6")
$rangeG.Interior.Color = 5296274

# Fill in the "Done!" status for rows 5 and 6 (status column G)
$ws.Range("G5").Value = "Done!"
$ws.Range("G6").Value = "Done!"

# Update row 10's Machine / GPU column: experiment re-run all on DLT1, now DLT1 / 3
$ws.Range("C10").Value = "DLT1 / 3"

# Set width for new column G (closest reachable value to the target 19.7109375)
$ws.Columns.Item(7).ColumnWidth = 18.8

# Leave selection on E5
$ws.Range("E5").Select()
